# The deck was reopened/resaved a day later than it was originally
# authored, so every auto-updating "datetimeFigureOut" date field
# (slide master, all slide layouts, and the notes master) needs its
# cached display text bumped from 1/10/2021 to 1/11/2021.

$p = $ppt.ActivePresentation

$oldDate = "1/10/2021"
$newDate = "1/11/2021"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own Date placeholder.
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout (CustomLayout) hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master's Date placeholder.
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
